$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C_32")

# Spanish month abbreviations, keyed by the month number (1-12) currently
# stored in column C ("Mes") of the Tabla3 table.
$months = @{
    1 = "Ene."
    2 = "Feb."
    3 = "Mar."
    4 = "Abr."
    5 = "May."
    6 = "Jun."
    7 = "Jul."
    8 = "Ago."
    9 = "Sep."
    10 = "Oct."
    11 = "Nov."
    12 = "Dic."
}

for ($r = 6; $r -le 85; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $num = [int]$cell.Value()
    $cell.Value = $months[$num]
}

